$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '287.19'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '1.49%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '29.58'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '4.11%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.124'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '1.56%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06704'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '3.22%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '7.340'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '1.57%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.396'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '1.04%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.381'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-1.48%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9207'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.33%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1598'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '4.03%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06824'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '5.22%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07765'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '2.00%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.02926'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '4.73%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.08983'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.17%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001573'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.81%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.04496'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '1.21%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0006476'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '1.88%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.006246'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '3.02%'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.02%'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-0.61%'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-2.83%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.088'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '2.44%'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '2.42%'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '0.81%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004125'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-7.51%'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-0.14%'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '-0.16%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04275'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '3.99%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006761'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '2.23%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1240'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '0.70%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002216'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '3.11%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01211'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '5.11%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005709'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '5.84%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.974'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '2.14%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-29.49%'
